$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1800193333333333
$ws.Range("H2").Value = 0.5400579999999999
$ws.Range("M2").Value = 71.44418333333333
$ws.Range("N2").Value = 214.33255
$ws.Range("O2").Value = 0.6986063918429039
$ws.Range("P2").Value = 0.6986063918429037
$ws.Range("Q2").Value = 12.86133425421111
$ws.Range("R2").Value = 115.7520082879
$ws.Range("S2").Value = 0.6986063918429039
$ws.Range("T2").Value = 0.6986063918429037

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1800193333333333
$ws.Range("H3").Value = 0.5400579999999999
$ws.Range("O3").Value = 0.1188372961583501
$ws.Range("P3").Value = 0.1188372961583501
$ws.Range("Q3").Value = 2.187793019939777
$ws.Range("R3").Value = 19.690137179458
$ws.Range("S3").Value = 0.1188372961583501
$ws.Range("T3").Value = 0.1188372961583501

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1800193333333333
$ws.Range("H4").Value = 0.5400579999999999
$ws.Range("M4").Value = 18.571964
$ws.Range("N4").Value = 55.715892
$ws.Range("O4").Value = 0.1816032062252276
$ws.Range("P4").Value = 0.1816032062252276
$ws.Range("Q4").Value = 3.343312577970666
$ws.Range("R4").Value = 30.08981320173599
$ws.Range("S4").Value = 0.1816032062252276
$ws.Range("T4").Value = 0.1816032062252276

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1800193333333333
$ws.Range("H5").Value = 0.5400579999999999
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09747100000000002
$ws.Range("N5").Value = 0.292413
$ws.Range("O5").Value = 0.000953105773518577
$ws.Range("P5").Value = 0.0009531057735185768
$ws.Range("Q5").Value = 0.01754666443933333
$ws.Range("R5").Value = 0.157919979954
$ws.Range("S5").Value = 0.000953105773518577
$ws.Range("T5").Value = 0.0009531057735185768
